# Adds the 2022-Q4 sheet (new quarterly snapshot) to the workbook, inserted
# right after "总计" and before "2022-Q3", and updates the "总计" summary
# sheet with a new top row for 2022-Q4 (shifting the existing rows down).
#
# All other existing quarter sheets (2022-Q3 .. 2021-Q1) keep their data and
# tab names unchanged; they simply end up one position further right because
# the new sheet is inserted before them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q4,
#    pushing the existing 6 rows down by one, and fill in the new row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Copy formatting of the last populated data row down into the new last row
# (row 8) so the newly-exposed row carries the same per-column styling as
# every other data row, then shift the B:D data (label/count/value) down by
# one row (working bottom-up so we never clobber a row before reading it).
# Column A is just a fixed 0-based row counter (0,1,2,...) that does NOT
# travel with the data - it is left as-is for rows 2-7 and only the new
# row 8 needs a value (6) appended.
$summary.Range("A7:D7").Copy()
$summary.Range("A8:D8").PasteSpecial(-4122)

for ($r = 7; $r -ge 2; $r--) {
    $dst = $r + 1
    $summary.Cells.Item($dst, 2).Value2 = $summary.Cells.Item($r, 2).Value2
    $summary.Cells.Item($dst, 3).Value2 = $summary.Cells.Item($r, 3).Value2
    $summary.Cells.Item($dst, 4).Value2 = $summary.Cells.Item($r, 4).Value2
}

$summary.Cells.Item(8, 1).Value2 = 6
$summary.Cells.Item(2, 2).Value2 = "2022-Q4"
$summary.Cells.Item(2, 3).Value2 = 5
$summary.Cells.Item(2, 4).Value2 = 0.54

# ---------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q4" sheet before the existing "2022-Q3"
#    sheet (i.e. right after "总计"), by duplicating "2022-Q3" (so it
#    inherits identical column widths / header styling / page setup)
#    and then overwriting its contents with the 2022-Q4 dataset.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Move it so it sits immediately before "2022-Q3" (it already does, since
# Copy($q3) places the new sheet right before $q3 - nothing further to do).

# Clear out the rows copied over from 2022-Q3 (rows 2-12) beyond what the
# 2022-Q4 dataset needs (rows 2-6), then fill in the new data.
$q4.Rows("7:12").Delete()

$q4Data = @(
    @(0, "210009", "金鹰核心资源混合", "3.14", "91.78", "5.00", "0.1570", 5),
    @(1, "001167", "金鹰科技创新股票", "3.17", "91.02", "4.72", "0.1496", 7),
    @(2, "162102", "金鹰中小盘精选混合", "3.48", "76.23", "4.29", "0.1493", 7),
    @(3, "210002", "金鹰红利价值混合A", "0.96", "77.22", "5.53", "0.0531", 5),
    @(4, "016563", "金鹰红利价值混合C", "0.52", "77.22", "5.53", "0.0288", 5)
)

# B (fund code, e.g. "001167" with a significant leading zero) and D:G
# (numeric-looking text, e.g. "3.14") hold text in the source data files,
# so force those columns to Text before writing them, otherwise they'd be
# auto-coerced to numbers (and leading zeros would be lost).
$q4.Range("B2:B6").NumberFormat = "@"
$q4.Range("D2:G6").NumberFormat = "@"

$rowIdx = 2
foreach ($row in $q4Data) {
    $q4.Cells.Item($rowIdx, 1).Value2 = $row[0]
    $q4.Cells.Item($rowIdx, 2).Value2 = $row[1]
    $q4.Cells.Item($rowIdx, 3).Value2 = $row[2]
    $q4.Cells.Item($rowIdx, 4).Value2 = $row[3]
    $q4.Cells.Item($rowIdx, 5).Value2 = $row[4]
    $q4.Cells.Item($rowIdx, 6).Value2 = $row[5]
    $q4.Cells.Item($rowIdx, 7).Value2 = $row[6]
    $q4.Cells.Item($rowIdx, 8).Value2 = $row[7]
    $rowIdx++
}

# ---------------------------------------------------------------------
# 3) Restore the "active sheet" tab to the last sheet (2021-Q1), matching
#    the original workbook where that sheet was the selected one.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
